# Updates the cryptocurrency price/volume snapshot values in the
# "cryptos" worksheet (Sheet1) to the latest scraped figures.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
# Some rows' relative ranking changed, so a handful of rows have their
# entire Coin/Link/Price/Volume values swapped with a neighboring row.
#
# Price values are stored as plain text in this sheet (they use a "."
# as a thousands separator, e.g. "66.077.83", which Excel would
# otherwise misread as a number). For the purely-numeric-looking price
# strings we prefix the literal with an apostrophe to force Excel to
# keep storing it as text, then reset the cell style back to "Normal"
# so the quote-prefix formatting flag doesn't linger on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.962.85'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '3.768.99'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = "'426.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.77%  '
$ws.Range("D6").Value = "'138.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.55%  '
$ws.Range("D7").Value = "'0.621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.70%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = "'0.727"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.74%  '
$ws.Range("D10").Value = "'0.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -10.74%  '
$ws.Range("D11").Value = "'0.0000303"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -14.00%  '
$ws.Range("D12").Value = "'42.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.45%  '
$ws.Range("D13").Value = "'10.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.98%  '
$ws.Range("D14").Value = '4.378.98'
$ws.Range("E14").Value = '  +0.31%  '
$ws.Range("D15").Value = "'15.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.53%  '
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = '3.746.58'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Value = "'19.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.00%  '
$ws.Range("E19").Value = '  +5.25%  '
$ws.Range("D20").Value = '66.108.42'
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = "'402.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.04%  '
$ws.Range("D22").Value = "'14.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.64%  '
$ws.Range("D23").Value = "'3.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.05%  '
$ws.Range("D24").Value = "'84.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.69%  '
$ws.Range("B25").Value = 'EthereumClassic'
$ws.Range("C25").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D25").Value = "'36.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.06%  '
$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").Value = "'9.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +34.85%  '
$ws.Range("D27").Value = "'3.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.50%  '
$ws.Range("D28").Value = "'9.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.34%  '
$ws.Range("E29").Value = '  -4.68%  '
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").Value = "'709.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.93%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = "'13.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +11.09%  '
$ws.Range("E32").Value = '  +11.99%  '
$ws.Range("E33").Value = '  +1.64%  '
$ws.Range("D34").Value = "'40.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.95%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").Value = "'5.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +34.85%  '
$ws.Range("E37").Value = '  -3.12%  '
$ws.Range("D38").Value = "'56.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.10%  '
$ws.Range("D39").Value = "'0.0468"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.67%  '
$ws.Range("D40").Value = "'2.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +36.72%  '
$ws.Range("E41").Value = '  +1.81%  '
$ws.Range("E42").Value = '  +4.47%  '
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("D44").Value = '0.0₃0657'
$ws.Range("E44").Value = '  -9.68%  '
$ws.Range("D45").Value = "'3.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.39%  '
$ws.Range("B46").Value = 'LidoDAOToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D46").Value = "'3.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.44%  '
$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D47").Value = "'0.321"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.55%  '
$ws.Range("D48").Value = "'2.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.27%  '
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("D50").Value = "'138.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'2.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.23%  '
